$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.070405244827271
$ws.Range("B1").Value = 2.523457288742065
$ws.Range("C1").Value = 2.618856430053711
$ws.Range("D1").Value = 3.253490447998047
$ws.Range("E1").Value = 0.8091923594474792
